# Update contact emails for the two remaining rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "priyank.mali.5057@gmail.com"
$ws.Range("B3").Value = "priyankmali297@gmail.com"

# Remove the extra sample rows (Alice Johnson ... Henry Lewis), keeping
# only the header row plus John Doe / Jane Smith.
$ws.Range("A4:A11").EntireRow.Delete()

# Adjust row heights to match the re-saved layout.
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 17.25
$ws.Rows.Item(3).RowHeight = 17.25
